$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": insert new daily row (2020-05-14 / serial 43965) before
# the trailing footnote row, pushing the footnote from row 37 to row 38.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()
$wsAll.Rows.Item(37).Insert()
$wsAll.Range("A37").Value = 43965
$wsAll.Range("B37").Value = 278
$wsAll.Range("C37").Value = 276
$wsAll.Range("D37").Value = 77
$wsAll.Range("E37").Value = 67
$wsAll.Range("F37").Value = 10
$wsAll.Range("G37").Value = 11
$wsAll.Range("H37").Value = 188
$wsAll.Range("I37").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "kobe": same pattern, footnote moves from row 92 to row 93.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()
$wsKobe.Rows.Item(92).Insert()
$wsKobe.Range("A92").Value = 43965
$wsKobe.Range("B92").Value = 0
$wsKobe.Range("C92").Value = 2743
$wsKobe.Range("D92").Value = 0
$wsKobe.Range("E92").Value = 281
$wsKobe.Range("F92").Value = 72
$wsKobe.Range("G92").Value = 63
$wsKobe.Range("H92").Value = 9
$wsKobe.Range("I92").Value = 11
$wsKobe.Range("J92").Value = 179
$wsKobe.Range("A92").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "other": same pattern, footnote moves from row 67 to row 68.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()
$wsOther.Rows.Item(67).Insert()
$wsOther.Range("A67").Value = 43965
$wsOther.Range("B67").Value = 0
$wsOther.Range("C67").Value = 14
$wsOther.Range("D67").Value = 5
$wsOther.Range("E67").Value = 4
$wsOther.Range("F67").Value = 1
$wsOther.Range("G67").Value = 0
$wsOther.Range("H67").Value = 9
$wsOther.Range("A67").Select() | Out-Null

# Restore original active sheet ("all" is tabSelected in the source file).
$wsAll.Activate()
